$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row (43) with the next month's data (01-07-2021)
# Force column A to be treated as text so the date-like label "01-07-2021"
# is stored as a string (matching the existing date-label cells) instead of
# being auto-converted into a date serial number.
$ws.Range("A43").NumberFormat = "@"
$ws.Range("A43").Value = "01-07-2021"
$ws.Range("A43").Style = "Normal"

$ws.Range("B43").Value = 0.8
$ws.Range("C43").Value = 0.6
$ws.Range("D43").Value = 1.2
$ws.Range("E43").Value = 0.4
$ws.Range("F43").Value = -0.5
